$wb = $excel.ActiveWorkbook
$new = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new.Range("E10").Value = "Layer"
$new.Range("F10").Value = "Smaller Model`n(N)"
$new.Range("G10").Value = "Large Model `n(N)"
$rng = $new.Range("E10:G10")
$rng.Font.Bold = $true
$rng.Font.Color = 16777215   # white
$rng.Interior.Color = 0      # black
$rng.HorizontalAlignment = -4108  # xlCenter
$new.Range("F10:G10").WrapText = $true
$rng.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$rng.Borders.Item(10).Weight = 2      # xlThin
$rng.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$rng.Borders.Item(8).Weight = 2
$new.Rows.Item(10).RowHeight = 34
